# Update cached market-price / profit snapshot values across all profession sheets.
# Values below are literal doubles (bit-identical to source data); Excel's own
# shortest round-trip float formatter may render them slightly differently on save,
# which is cosmetic only (same IEEE-754 value).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J17").Value2 = 564.7778
$ws.Range("N17").Value2 = -2030.3334
$ws.Range("L17").Value2 = 1694.3334
$ws.Range("H17").Value2 = 577.1053
$ws.Range("K41").Value2 = 157.66667
$ws.Range("H41").Value2 = 334.4
$ws.Range("M41").Value2 = 282.33333
$ws.Range("I41").Value2 = 157.66667
$ws.Range("I86").Value2 = 5638.0
$ws.Range("H86").Value2 = 6460.1816
$ws.Range("M86").Value2 = -4515.0
$ws.Range("K86").Value2 = 5638.0
$ws.Range("M88").Value2 = -11112044.0
$ws.Range("H88").Value2 = 4001936.8
$ws.Range("L88").Value2 = 2272.9375
$ws.Range("I88").Value2 = 11112450.0
$ws.Range("J88").Value2 = 2272.9375
$ws.Range("N88").Value2 = -3084.9375
$ws.Range("K88").Value2 = 11112450.0
$ws.Range("M89").Value2 = -22574.0
$ws.Range("K89").Value2 = 28190.0
$ws.Range("H89").Value2 = 6460.1816
$ws.Range("I89").Value2 = 5638.0
$ws.Range("M91").Value2 = -11111046.0
$ws.Range("N91").Value2 = -5080.9375
$ws.Range("J91").Value2 = 2272.9375
$ws.Range("L91").Value2 = 2272.9375
$ws.Range("H91").Value2 = 4001936.8
$ws.Range("K91").Value2 = 11112450.0
$ws.Range("I91").Value2 = 11112450.0
$ws.Range("J92").Value2 = 659.5
$ws.Range("K92").Value2 = 981.125
$ws.Range("L92").Value2 = 659.5
$ws.Range("I92").Value2 = 981.125
$ws.Range("N92").Value2 = -3155.5
$ws.Range("H92").Value2 = 916.8
$ws.Range("M92").Value2 = 266.875
$ws.Range("N113").Value2 = -12185.0
$ws.Range("L113").Value2 = 5677.0
$ws.Range("J113").Value2 = 5677.0
$ws.Range("M132").Value2 = -9918.215900000001
$ws.Range("K132").Value2 = 12448.2159
$ws.Range("H132").Value2 = 4007.1025
$ws.Range("I132").Value2 = 4149.4053
$ws.Range("M141").Value2 = -4442.125
$ws.Range("K141").Value2 = 9622.125
$ws.Range("H141").Value2 = 3406.4443
$ws.Range("I141").Value2 = 3207.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1286.7142
$ws.Range("K2").Value2 = 1286.7142
$ws.Range("I2").Value2 = 1286.7142
$ws.Range("M2").Value2 = -1173.7142
$ws.Range("K116").Value2 = 1286.7142
$ws.Range("M116").Value2 = 1007.2858
$ws.Range("H116").Value2 = 1286.7142
$ws.Range("I116").Value2 = 1286.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1286.7142
$ws.Range("K3").Value2 = 1286.7142
$ws.Range("I3").Value2 = 1286.7142
$ws.Range("M3").Value2 = -1172.7142
$ws.Range("K80").Value2 = 1000000000.0
$ws.Range("H80").Value2 = 200000350.0
$ws.Range("J80").Value2 = 446.75
$ws.Range("I80").Value2 = 1000000000.0
$ws.Range("M80").Value2 = -999999002.0
$ws.Range("N80").Value2 = -2442.75
$ws.Range("L80").Value2 = 446.75
$ws.Range("M83").Value2 = -4999995008.0
$ws.Range("K83").Value2 = 5000000000.0
$ws.Range("N83").Value2 = -12217.75
$ws.Range("I83").Value2 = 1000000000.0
$ws.Range("J83").Value2 = 446.75
$ws.Range("L83").Value2 = 2233.75
$ws.Range("H83").Value2 = 200000350.0
$ws.Range("H99").Value2 = 2215.6
$ws.Range("M99").Value2 = 173.9000000000001
$ws.Range("K99").Value2 = 1324.1
$ws.Range("I99").Value2 = 1324.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 5684416.0
$ws.Range("I31").Value2 = 1608.2354
$ws.Range("K31").Value2 = 1608.2354
$ws.Range("M31").Value2 = -1313.2354
$ws.Range("I34").Value2 = 1608.2354
$ws.Range("H34").Value2 = 5684416.0
$ws.Range("K34").Value2 = 1608.2354
$ws.Range("M34").Value2 = -1406.2354
$ws.Range("J99").Value2 = 4385.8887
$ws.Range("H99").Value2 = 4012.6
$ws.Range("M99").Value2 = -1954.6667
$ws.Range("N99").Value2 = -7381.8887
$ws.Range("L99").Value2 = 4385.8887
$ws.Range("K99").Value2 = 3452.6667
$ws.Range("I99").Value2 = 3452.6667
$ws.Range("H122").Value2 = 3245.625
$ws.Range("N122").Value2 = -18467.9995
$ws.Range("L122").Value2 = 13567.9995
$ws.Range("J122").Value2 = 4522.6665
$ws.Range("H126").Value2 = 4012.6
$ws.Range("L126").Value2 = 13157.6661
$ws.Range("M126").Value2 = -7888.000100000001
$ws.Range("I126").Value2 = 3452.6667
$ws.Range("K126").Value2 = 10358.0001
$ws.Range("J126").Value2 = 4385.8887
$ws.Range("N126").Value2 = -18097.6661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value2 = 1787822.0
$ws.Range("K4").Value2 = 4222540.800000001
$ws.Range("I4").Value2 = 1407513.6
$ws.Range("M4").Value2 = -4222428.800000001
$ws.Range("M113").Value2 = -599.0
$ws.Range("I113").Value2 = 923.0
$ws.Range("N113").Value2 = -7934.857400000001
$ws.Range("K113").Value2 = 2769.0
$ws.Range("L113").Value2 = 3594.8574
$ws.Range("J113").Value2 = 1198.2858
$ws.Range("H113").Value2 = 1115.7
$ws.Range("N131").Value2 = -11780627.4
$ws.Range("J131").Value2 = 3923515.8
$ws.Range("L131").Value2 = 11770547.4
$ws.Range("H131").Value2 = 2567641.8
$ws.Range("H134").Value2 = 3896.1333
$ws.Range("K134").Value2 = 8790.75
$ws.Range("M134").Value2 = -3720.75
$ws.Range("I134").Value2 = 2930.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J109").Value2 = 120000.0
$ws.Range("L109").Value2 = 120000.0
$ws.Range("N109").Value2 = -122080.0
$ws.Range("H109").Value2 = 120000.0
$ws.Range("M132").Value2 = -3841.8572
$ws.Range("K132").Value2 = 6371.8572
$ws.Range("H132").Value2 = 2184.3635
$ws.Range("I132").Value2 = 2123.9524

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 1024.7142
$ws.Range("L16").Value2 = 0.0
$ws.Range("I16").Value2 = 1024.7142
$ws.Range("K16").Value2 = 1024.7142
$ws.Range("J16").Value2 = 0.0
$ws.Range("N16").ClearContents()
$ws.Range("M16").Value2 = -854.7141999999999
$ws.Range("J22").Value2 = 2511.0
$ws.Range("K22").Value2 = 2352.75
$ws.Range("L22").Value2 = 2511.0
$ws.Range("H22").Value2 = 2384.4
$ws.Range("M22").Value2 = -2057.75
$ws.Range("I22").Value2 = 2352.75
$ws.Range("N22").Value2 = -3101.0
$ws.Range("I27").Value2 = 2352.75
$ws.Range("L27").Value2 = 2511.0
$ws.Range("J27").Value2 = 2511.0
$ws.Range("N27").Value2 = -2725.0
$ws.Range("K27").Value2 = 2352.75
$ws.Range("M27").Value2 = -2245.75
$ws.Range("H27").Value2 = 2384.4
$ws.Range("K61").Value2 = 1398.25
$ws.Range("M61").Value2 = -1196.25
$ws.Range("I61").Value2 = 1398.25
$ws.Range("H61").Value2 = 1398.25
$ws.Range("I100").Value2 = 2184.2856
$ws.Range("J100").Value2 = 1361.0
$ws.Range("H100").Value2 = 1937.3
$ws.Range("K100").Value2 = 2184.2856
$ws.Range("M100").Value2 = -1643.2856
$ws.Range("N100").Value2 = -2443.0
$ws.Range("L100").Value2 = 1361.0
$ws.Range("M113").Value2 = 771.75
$ws.Range("I113").Value2 = 1398.25
$ws.Range("K113").Value2 = 1398.25
$ws.Range("H113").Value2 = 1398.25
$ws.Range("M132").Value2 = -9164.75
$ws.Range("K132").Value2 = 11694.75
$ws.Range("H132").Value2 = 3966.0
$ws.Range("L132").Value2 = 12168.9999
$ws.Range("N132").Value2 = -17228.9999
$ws.Range("I132").Value2 = 3898.25
$ws.Range("J132").Value2 = 4056.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value2 = 500000.0
$ws.Range("L135").Value2 = 500000.0
$ws.Range("N135").Value2 = -510140.0
$ws.Range("J135").Value2 = 500000.0

